$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.5149680659622189
$ws.Range("C2").Value = 0.5621490832682646
$ws.Range("D2").Value = 0.398099301564984
$ws.Range("E2").Value = 0.6309511086962158
$ws.Range("F2").Value = 0.3773599919233752
$ws.Range("G2").Value = 15

$ws.Range("B3").Value = 0.374392881790439
$ws.Range("C3").Value = 0.3867115873339156
$ws.Range("D3").Value = 0.2095230518030864
$ws.Range("E3").Value = 0.4577368805362819
$ws.Range("F3").Value = 0.2732908152567819
$ws.Range("G3").Value = 14

$ws.Range("B4").Value = 0.3047392045302661
$ws.Range("C4").Value = 0.3200881600767664
$ws.Range("D4").Value = 0.1574586298430795
$ws.Range("E4").Value = 0.3968105717380517
$ws.Range("F4").Value = 0.2645285762523938
$ws.Range("G4").Value = 13

$ws.Range("B5").Value = 0.4461589838474358
$ws.Range("C5").Value = 0.4461589838474358
$ws.Range("D5").Value = 0.2557785926581425
$ws.Range("E5").Value = 0.5057455809576021
$ws.Range("F5").Value = 0.248751253168957
$ws.Range("G5").Value = 12

$ws.Range("B6").Value = 0.410700336148952
$ws.Range("C6").Value = 0.4135187035930826
$ws.Range("D6").Value = 0.2265423492358624
$ws.Range("E6").Value = 0.4759646512461428
$ws.Range("F6").Value = 0.2522981201580784
$ws.Range("G6").Value = 11

$ws.Range("B7").Value = 0.3494740007469088
$ws.Range("C7").Value = 0.3597227794264332
$ws.Range("D7").Value = 0.1638948727701445
$ws.Range("E7").Value = 0.4048393172236913
$ws.Range("F7").Value = 0.2154138022300699
$ws.Range("G7").Value = 10

$ws.Range("B8").Value = 0.3706545189354704
$ws.Range("C8").Value = 0.3734090874777368
$ws.Range("D8").Value = 0.1807682845776609
$ws.Range("E8").Value = 0.4251685366741769
$ws.Range("F8").Value = 0.2209218214474817
$ws.Range("G8").Value = 9

$ws.Range("B9").Value = 0.3834281208488316
$ws.Range("C9").Value = 0.385141425687041
$ws.Range("D9").Value = 0.1960677877649145
$ws.Range("E9").Value = 0.4427954242818171
$ws.Range("F9").Value = 0.2367654991933658
$ws.Range("G9").Value = 8

$ws.Range("B10").Value = 0.3533129028593847
$ws.Range("C10").Value = 0.3533129028593847
$ws.Range("D10").Value = 0.1644395682562605
$ws.Range("E10").Value = 0.4055114896723155
$ws.Range("F10").Value = 0.2149677985595936
$ws.Range("G10").Value = 7

$ws.Range("B11").Value = 0.3494951052803065
$ws.Range("C11").Value = 0.3494951052803065
$ws.Range("D11").Value = 0.1636482358169425
$ws.Range("E11").Value = 0.404534591619731
$ws.Range("F11").Value = 0.2231629195060417
$ws.Range("G11").Value = 6

